# Apply the "glossary sort + tab-panel label" edit described by the commit.
$wb = $excel.ActiveWorkbook

# --- 1. Sort the glossary table (term, definition, source) A-Z by term ---
$glossary = $wb.Worksheets.Item("glossary")
$sortRange = $glossary.Range("A2:C20")
$keyRange  = $glossary.Range("A2:A20")

$glossary.Sort.SortFields.Clear()
$glossary.Sort.SortFields.Add($keyRange, 0, 1, 0, 0)
$glossary.Sort.SetRange($sortRange)
$glossary.Sort.Header = 0
$glossary.Sort.MatchCase = $false
$glossary.Sort.Orientation = 1
$glossary.Sort.Apply()

# --- 2. Adjust tab panel label styles (overview_GS & overview_IS) ---
$overviewGS = $wb.Worksheets.Item("overview_GS")
$overviewGS.Range("B2").Interior.Pattern = -4142
$overviewGS.Range("B2").WrapText = $true

$overviewIS = $wb.Worksheets.Item("overview_IS")
$overviewIS.Range("B2").Interior.Pattern = -4142
$overviewIS.Range("B2").WrapText = $true

# --- 3. Make glossary sheet the active sheet/tab ---
$glossary.Activate()
$glossary.Range("A2:XFD20").Select()

$wb.Save()
